# Fruta / hortaliza, semanal
# Insert a new weekly record at row 204 (pushing all existing records for
# rows 204..270 down by one, to 205..271) and populate the new row with the
# latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 204:270 down to 205:271, leaving a blank row 204 to fill in.
$ws.Rows.Item(204).Insert()

$ws.Range("A204").Value = 3
$ws.Range("B204").Value = "Femacal de La Calera"
$ws.Range("C204").Value = "Coquimbo"
$ws.Range("D204").Value = 45119
$ws.Range("E204").Value = 5
$ws.Range("F204").Value = 100112026
$ws.Range("G204").Value = "Haba"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 40
$ws.Range("K204").Value = 17000
$ws.Range("L204").Value = 17000
$ws.Range("M204").Value = 17000
$ws.Range("N204").Value = "$/saco 25 kilos"
$ws.Range("O204").Value = "Provincia de Limarí"
$ws.Range("P204").Value = 680
$ws.Range("Q204").Value = 25
$ws.Range("R204").Value = "Hortaliza"
